# =====================================================================
# Weekly CompStat update (121st Precinct) - new crime data collected
# Updates the report header (volume number + week-covering dates) and
# refreshes the weekly/28-day/YTD/2-year crime-count figures in rows
# 16-30 of the CompStat grid.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 30   Number  9" -> "Volume 30   Number  10"
# (rich-text cell A8; replace only the trailing volume-number run so
#  the rest of the run-formatted string is left untouched)
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 1).Text = "10"

# ---------------------------------------------------------------------
# Header: "Report Covering the Week  2/27/2023  Through  3/5/2023"
#      -> "Report Covering the Week  3/6/2023  Through  3/12/2023"
# (rich-text cell C9; edit the right-most date run first so the left
#  run's character offsets are not shifted by the length change)
# ---------------------------------------------------------------------
$ws.Range("C9").Characters(47, 8).Text = "3/12/2023"
$ws.Range("C9").Characters(27, 9).Text = "3/6/2023"

# ---------------------------------------------------------------------
# Crime-statistics grid (rows 16-30): refreshed weekly counts and the
# derived percent-change figures.
# ---------------------------------------------------------------------

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -57.142857142857
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = -6.666666666666
$ws.Range("L16").Value = 55.555555555555
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 6.666666666666
$ws.Range("I17").Value = 45
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 28.571428571428
$ws.Range("L17").Value = 73.076923076923
$ws.Range("D18").Value = 4
$ws.Range("G18").Value = 12
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = -70
$ws.Range("L18").Value = -70
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -56.410256410256
$ws.Range("I19").Value = 51
$ws.Range("J19").Value = 91
$ws.Range("K19").Value = -43.956043956044
$ws.Range("L19").Value = -5.555555555555
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -17.647058823529
$ws.Range("L20").Value = 55.555555555555
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -43.421052631578
$ws.Range("I21").Value = 132
$ws.Range("J21").Value = 180
$ws.Range("K21").Value = -26.666666666666
$ws.Range("L21").Value = 9.090909090909
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -4.545454545454
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = 1.176470588235
$ws.Range("I24").Value = 233
$ws.Range("J24").Value = 222
$ws.Range("K24").Value = 4.954954954954
$ws.Range("L24").Value = 30.167597765363
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 18.75
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 81
$ws.Range("K25").Value = 19.753086419753
$ws.Range("L25").Value = 64.406779661017
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 13
$ws.Range("K27").Value = 44.444444444444
$ws.Range("L27").Value = 62.5

# --- Precincts that previously showed no activity now report real
#     figures: swap the placeholder-text format for the normal
#     numeric format used elsewhere in the column, then write the
#     figure. ---

$ws.Range("F23").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 2
$ws.Range("F23").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("K15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 100
$ws.Range("F23").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("K15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("F23").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("K15").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("F23").Copy()
$ws.Range("J30").PasteSpecial(-4122)
$ws.Range("J30").Value = 1
$ws.Range("K15").Copy()
$ws.Range("K30").PasteSpecial(-4122)
$ws.Range("K30").Value = -100

# --- Precincts that now have no activity switch back to the "0" /
#     "***.*" placeholder text (and its associated style) used for
#     the other zero-activity cells in the grid. ---

$ws.Range("D27").Formula = "'0"
$ws.Range("C23").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Formula = "'***.*"
$ws.Range("C23").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D28").Formula = "'0"
$ws.Range("C23").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Formula = "'***.*"
$ws.Range("C23").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").Formula = "'0"
$ws.Range("C23").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Formula = "'***.*"
$ws.Range("C23").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
